# Estadisticos Segundo Parcial 23 Mayo
# On sheet "Rescatables", rows 2 and 3 (the two rescatable students) swap
# their identity data (NC, Paterno, Materno, Nombres, Nombre_Largo, Grupo)
# while the "Reprobadas" count (column G) stays associated with the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Capture current (before) values for columns A-F on rows 2 and 3
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$d2 = $ws.Range("D2").Value()
$e2 = $ws.Range("E2").Value()
$f2 = $ws.Range("F2").Value()

$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()
$c3 = $ws.Range("C3").Value()
$d3 = $ws.Range("D3").Value()
$e3 = $ws.Range("E3").Value()
$f3 = $ws.Range("F3").Value()

# Write row 3's identity data into row 2, and row 2's identity data into row 3.
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3
$ws.Range("E2").Value = $e3
$ws.Range("F2").Value = $f3

$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2
$ws.Range("E3").Value = $e2
$ws.Range("F3").Value = $f2
